$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateUpdates = @{
    "B2" = "05/17/2023"
    "B3" = "05/24/2023"
    "B4" = "05/27/2023"
    "B5" = "05/16/2023"
    "B6" = "05/22/2023"
    "B7" = "05/17/2023"
    "B8" = "05/26/2023"
}

foreach ($addr in $dateUpdates.Keys) {
    $cell = $ws.Range($addr)
    # Force text entry so Excel doesn't reinterpret the MM/DD/YYYY string
    # as a date serial number - the source cells are plain text (inlineStr).
    $cell.NumberFormat = "@"
    $cell.Value = $dateUpdates[$addr]
    # Drop back to the default "Normal" style so no stray formatting is
    # left behind on the cell (matches the original unstyled cells).
    $cell.Style = "Normal"
}
